$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 11
# from 45185 (2023-09-16) to 45204 (2023-10-05)
$ws.Range("C2:C11").Value = 45204
